$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 29999
$ws.Range("I7").Value = 29999
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 29999
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -29887
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 29999
$ws.Range("I14").Value = 29999
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 29999
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -29808
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 617.375
$ws.Range("I32").Value = 900
$ws.Range("J32").Value = 397.55554
$ws.Range("K32").Value = 900
$ws.Range("L32").Value = 397.55554
$ws.Range("M32").Value = -574
$ws.Range("N32").Value = -1049.55554

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 12820989
$ws.Range("I33").Value = 22222768
$ws.Range("J33").Value = 380.9091
$ws.Range("K33").Value = 22222768
$ws.Range("L33").Value = 380.9091
$ws.Range("M33").Value = -22222539
$ws.Range("N33").Value = -838.9091000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2708.3333
$ws.Range("I64").Value = 2600
$ws.Range("J64").Value = 2730
$ws.Range("K64").Value = 2600
$ws.Range("L64").Value = 2730
$ws.Range("M64").Value = -2352
$ws.Range("N64").Value = -3226

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2708.3333
$ws.Range("I67").Value = 2600
$ws.Range("J67").Value = 2730
$ws.Range("K67").Value = 2600
$ws.Range("L67").Value = 2730
$ws.Range("M67").Value = -1742
$ws.Range("N67").Value = -4446

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 61564.883
$ws.Range("I76").Value = 65037.688
$ws.Range("J76").Value = 6000
$ws.Range("K76").Value = 65037.688
$ws.Range("L76").Value = 6000
$ws.Range("M76").Value = -64722.688
$ws.Range("N76").Value = -6630

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 61564.883
$ws.Range("I79").Value = 65037.688
$ws.Range("J79").Value = 6000
$ws.Range("K79").Value = 65037.688
$ws.Range("L79").Value = 6000
$ws.Range("M79").Value = -63945.688
$ws.Range("N79").Value = -8184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7333.3335
$ws.Range("J116").Value = 10000
$ws.Range("L116").Value = 10000
$ws.Range("N116").Value = -16884

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2133922.5
$ws.Range("I132").Value = 2343526.5
$ws.Range("J132").Value = 2949.3333
$ws.Range("K132").Value = 7030579.5
$ws.Range("L132").Value = 8847.999899999999
$ws.Range("M132").Value = -7028049.5
$ws.Range("N132").Value = -13907.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2885
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15874020
$ws.Range("I45").Value = 33334238
$ws.Range("K45").Value = 33334238
$ws.Range("M45").Value = -33333861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 10262.6
$ws.Range("I88").Value = 10903
$ws.Range("J88").Value = 9835.666999999999
$ws.Range("K88").Value = 10903
$ws.Range("L88").Value = 9835.666999999999
$ws.Range("M88").Value = -10497
$ws.Range("N88").Value = -10647.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 10262.6
$ws.Range("I91").Value = 10903
$ws.Range("J91").Value = 9835.666999999999
$ws.Range("K91").Value = 10903
$ws.Range("L91").Value = 9835.666999999999
$ws.Range("M91").Value = -9499
$ws.Range("N91").Value = -12643.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5437.6177
$ws.Range("I20").Value = 5768.04
$ws.Range("K20").Value = 5768.04
$ws.Range("M20").Value = -5521.04

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2585.2222
$ws.Range("I31").Value = 2347.1936
$ws.Range("J31").Value = 4061
$ws.Range("K31").Value = 2347.1936
$ws.Range("L31").Value = 4061
$ws.Range("M31").Value = -2052.1936
$ws.Range("N31").Value = -4651

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2585.2222
$ws.Range("I34").Value = 2347.1936
$ws.Range("J34").Value = 4061
$ws.Range("K34").Value = 2347.1936
$ws.Range("L34").Value = 4061
$ws.Range("M34").Value = -2145.1936
$ws.Range("N34").Value = -4465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 24449.75
$ws.Range("J74").Value = 29266.334
$ws.Range("L74").Value = 29266.334
$ws.Range("N74").Value = -31014.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 24449.75
$ws.Range("J77").Value = 29266.334
$ws.Range("L77").Value = 87799.00199999999
$ws.Range("N77").Value = -96535.00199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1004.44446
$ws.Range("J122").Value = 970
$ws.Range("L122").Value = 2910
$ws.Range("N122").Value = -7810

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 12503140
$ws.Range("I132").Value = 2004
$ws.Range("J132").Value = 20837230
$ws.Range("K132").Value = 6012
$ws.Range("L132").Value = 62511690
$ws.Range("M132").Value = -3482
$ws.Range("N132").Value = -62516750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5300
$ws.Range("I134").Value = 5504.5454
$ws.Range("K134").Value = 16513.6362
$ws.Range("M134").Value = -13978.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 22352.5
$ws.Range("I120").Value = 11598
$ws.Range("J120").Value = 25515.588
$ws.Range("K120").Value = 34794
$ws.Range("L120").Value = 76546.764
$ws.Range("M120").Value = -29956
$ws.Range("N120").Value = -86222.764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 30616344
$ws.Range("I70").Value = 40544604
$ws.Range("J70").Value = 4215.75
$ws.Range("K70").Value = 40544604
$ws.Range("L70").Value = 4215.75
$ws.Range("M70").Value = -40544334
$ws.Range("N70").Value = -4755.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 30616344
$ws.Range("I73").Value = 40544604
$ws.Range("J73").Value = 4215.75
$ws.Range("K73").Value = 40544604
$ws.Range("L73").Value = 4215.75
$ws.Range("M73").Value = -40543668
$ws.Range("N73").Value = -6087.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4782.773
$ws.Range("I80").Value = 4821.5
$ws.Range("J80").Value = 4750.5
$ws.Range("K80").Value = 4821.5
$ws.Range("L80").Value = 4750.5
$ws.Range("M80").Value = -3823.5
$ws.Range("N80").Value = -6746.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4782.773
$ws.Range("I83").Value = 4821.5
$ws.Range("J83").Value = 4750.5
$ws.Range("K83").Value = 24107.5
$ws.Range("L83").Value = 23752.5
$ws.Range("M83").Value = -19115.5
$ws.Range("N83").Value = -33736.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1892.2222
$ws.Range("I97").Value = 1608.75
$ws.Range("J97").Value = 4160
$ws.Range("K97").Value = 1608.75
$ws.Range("L97").Value = 4160
$ws.Range("M97").Value = -1112.75
$ws.Range("N97").Value = -5152

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1413.2759
$ws.Range("I102").Value = 1095.85
$ws.Range("J102").Value = 2118.6667
$ws.Range("K102").Value = 1095.85
$ws.Range("L102").Value = 2118.6667
$ws.Range("M102").Value = 526.1500000000001
$ws.Range("N102").Value = -5362.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 22728582
$ws.Range("I113").Value = 125000600
$ws.Range("K113").Value = 125000600
$ws.Range("M113").Value = -124998430

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4445.4
$ws.Range("I122").Value = 6478.1113
$ws.Range("J122").Value = 2782.2727
$ws.Range("K122").Value = 19434.3339
$ws.Range("L122").Value = 8346.8181
$ws.Range("M122").Value = -16984.3339
$ws.Range("N122").Value = -13246.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 38945
$ws.Range("J123").Value = 38945
$ws.Range("L123").Value = 38945
$ws.Range("N123").Value = -48745

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1897.0667
$ws.Range("I122").Value = 1833.0435
$ws.Range("J122").Value = 2107.4285
$ws.Range("K122").Value = 5499.1305
$ws.Range("L122").Value = 6322.2855
$ws.Range("M122").Value = -3049.1305
$ws.Range("N122").Value = -11222.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 32264236
$ws.Range("I126").Value = 45461668
$ws.Range("K126").Value = 136385004
$ws.Range("M126").Value = -136382534
